# Refresh the "Coin" symbol list / price feed, per the Dec 18 2022 GitHub
# Actions scrape. A handful of rows near the bottom of the ranked list
# swapped places with their neighbour (exchange tokens jockeying for
# position), and every "Price" cell got a freshly polled quote.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells are stored as literal text (e.g. "0.001660", "3.208") even
# though they look numeric - plain `.Value = "3.208"` would get coerced to
# the number 3.208 and lose the trailing zero. Forcing text via a leading
# apostrophe (classic Excel "treat as text" trick) and then resetting the
# cell style keeps the stored value an exact string match with no stray
# number formatting left behind.
function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $Sheet.Range($Addr).Value = "'" + $Val
    $Sheet.Range($Addr).Style = "Normal"
}

Set-TextValue $ws "D2" "246.18"
Set-TextValue $ws "D3" "22.43"
Set-TextValue $ws "D4" "5.516"
Set-TextValue $ws "D5" "0.05637"
Set-TextValue $ws "D6" "6.468"
Set-TextValue $ws "D7" "0.8053"
Set-TextValue $ws "D8" "1.054"
Set-TextValue $ws "D9" "0.1438"
Set-TextValue $ws "D10" "0.07391"
Set-TextValue $ws "D11" "0.03193"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D12" "0.02924"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D13" "0.09251"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D14" "0.001660"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D15" "3.208"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D16" "0.04728"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D17" "0.0005839"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D18" "0.006274"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D19" "0.001063"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D20" "0.004112"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D21" "0.0001506"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D22" "3.977"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "GateToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D23" "3.380"
$ws.Range("E23").Value = "22GateTokenGT"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D24" "2.140"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws "D25" "0.3274"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws "D26" "0.1312"
$ws.Range("E26").Value = "25ProBitTokenPROBBestin24h"
Set-TextValue $ws "D27" "0.0003004"
Set-TextValue $ws "D40" "0.04157"
Set-TextValue $ws "D41" "0.006889"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1041"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002983"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws "D44" "0.009842"
Set-TextValue $ws "D45" "0.00005645"
Set-TextValue $ws "D47" "0.6810"
Set-TextValue $ws "D48" "0.02079"
Set-TextValue $ws "D49" "0.00002103"
Set-TextValue $ws "D50" "0.01011"
